# Add two new readings to the bottom of the list, matching the green/no-fill
# "currently reading" look used by the other not-yet-highlighted rows
# (e.g. A17:A19), then leave the selection on the newest entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "Lutz 2019, Part VI, Classes and OOP"
$ws.Range("A17").Copy()
$ws.Range("A22").PasteSpecial(-4122)

$ws.Range("A23").Value = "Lutz, 2019, Part VII: Exceptions and Tools"
$ws.Range("A18").Copy()
$ws.Range("A23").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("A23").Select()
